$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "GNG_TO-16502911953442924"
$ws.Range("B2").Value = "go_stims-16502911952932935.csv"
$ws.Range("B3").Value = "GNG_stims-16502911953112948.csv"
$ws.Range("B4").Value = "go_stims-1650291195314294.csv"
$ws.Range("B5").Value = "GNG_stims-16502911953422925.csv"

# --- Sheet 2: NB_TO ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "NB_TO-16502911979552925"
$ws.Range("B2").Value = "ZB-match_4-1650291195438295.csv"
$ws.Range("B3").Value = "TB-16502911973922935.csv"
$ws.Range("B4").Value = "OB-16502911960902977.csv"
$ws.Range("B5").Value = "ZB-match_6-16502911955322938.csv"
$ws.Range("B6").Value = "OB-1650291195568293.csv"
$ws.Range("B7").Value = "ZB-match_1-16502911954552908.csv"
$ws.Range("B8").Value = "TB-1650291197794293.csv"
$ws.Range("B9").Value = "OB-16502911959142945.csv"
$ws.Range("B10").Value = "TB-16502911979323013.csv"

# --- Sheet 3: RS_TO ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "RS_TO-16502911979583113"
$ws.Range("B2").Value = "eyes closed"
$ws.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "TOL_TO-16502911980362954"
$ws.Range("B2").Value = "MM_stims-16502911979863157.csv"
$ws.Range("B3").Value = "ZM_stims-16502911979612958.csv"
$ws.Range("B4").Value = "MM_stims-16502911980172942.csv"
$ws.Range("B5").Value = "ZM_stims-16502911979872973.csv"
$ws.Range("B6").Value = "MM_stims-16502911980343003.csv"
$ws.Range("B7").Value = "ZM_stims-16502911980182931.csv"

# --- Sheet 5: vSAT_TO ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "vSAT_TO-16502911981112912"
$ws.Range("B2").Value = "vSAT_stims-16502911980952954.csv"
$ws.Range("B3").Value = "SAT_stims-16502911980652945.csv"
$ws.Range("B4").Value = "vSAT_stims-16502911980812948.csv"
$ws.Range("B5").Value = "SAT_stims-16502911980402954.csv"
